$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 749, pushing existing rows 749..809 down to 750..810
$ws.Rows.Item(749).Insert()

# Populate the newly inserted row 749 with the new record's data
$ws.Cells.Item(749, 1).Value = 6
$ws.Cells.Item(749, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(749, 3).Value = "Metropolitana"
$ws.Cells.Item(749, 4).Value = "2023-10-24"
$ws.Cells.Item(749, 5).Value = 13
$ws.Cells.Item(749, 6).Value = 100112052
$ws.Cells.Item(749, 7).Value = "Albahaca"
$ws.Cells.Item(749, 8).Value = "Sin especificar"
$ws.Cells.Item(749, 9).Value = "Primera"
$ws.Cells.Item(749, 10).Value = 250
$ws.Cells.Item(749, 11).Value = 5000
$ws.Cells.Item(749, 12).Value = 5000
$ws.Cells.Item(749, 13).Value = 5000
$ws.Cells.Item(749, 14).Value = "`$/docena de matas"
$ws.Cells.Item(749, 15).Value = "Región Metropolitana"
$ws.Cells.Item(749, 16).Value = 833
$ws.Cells.Item(749, 17).Value = 6
$ws.Cells.Item(749, 18).Value = "Hortaliza"
